$wb = $excel.ActiveWorkbook

# 1. Rename sheet "set_document" to "add_document"
$docSheet = $wb.Worksheets.Item("set_document")
$docSheet.Name = "add_document"

# Update the function prototype text on the renamed sheet.
$docSheet.Range("B2").Value = "function module:add_document(name, filename)"

# 2. Insert a new row in the "summary" sheet before row 17
#    for the new ":add_document" function entry.
$summary = $wb.Worksheets.Item("summary")
$summary.Rows.Item(17).Insert()

$summary.Range("A17").Value = ":add_document"
$summary.Range("B17").Value = "function"
$summary.Range("B17").HorizontalAlignment = -4108  # xlCenter
$summary.Range("C17").Value = "Add a document."

# 3. Update selections / active sheet state
$docSheet.Range("B4").Select()
$summary.Activate()
$summary.Range("A17").Select()
